$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "anasA" / bookmark "_GoBack" / "gA" runs into a single run
#    "anasAgA" (this also drops the now-stale _GoBack bookmark at that spot).
# ---------------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("anasA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStart = $f1.Start

$f2 = $d.Content
$f2.Find.Execute("gA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeEnd = $f2.End

$mergeRange = $d.Range($mergeStart, $mergeEnd)
# Force a real content change (same-text assignment is a no-op in this
# engine) so the bookmark actually collapses, then restore the final text.
$mergeRange.Text = "__TMP_PLACEHOLDER__"
$mergeRange2 = $d.Range($mergeStart, $mergeStart + 20)
$mergeRange2.Text = "anasAgA"

# ---------------------------------------------------------------------------
# 2) Append a new test-case row (TC_003 / logout) to the table.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()

# Cell 1: Test Case ID
$cell1 = $newRow.Cells.Item(1)
$cell1.Range.Text = "TC_003"

# Cell 2: Title
$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "Sign out from account and  return to the homepage."

# Cell 3: Description (two paragraphs, with a superscript "nd")
$cell3 = $newRow.Cells.Item(3)
$cell3.Range.Text = "(1) Execute 2nd test case.`r(2) Hit logout button."

$ndFind = $d.Content
$ndFind.Find.Execute("nd test case", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ndRange = $d.Range($ndFind.Start, $ndFind.Start + 2)
$ndRange.Font.Superscript = $true

# Cell 4: Pre-condition
$cell4 = $newRow.Cells.Item(4)
$cell4.Range.Text = "TC_002"

# Cell 5: Test Data -> left empty

# Cell 6: Expected result
$cell6 = $newRow.Cells.Item(6)
$cell6.Range.Text = "It should return back to the homepage after successful logout."

# Cell 7: Actual result
$cell7 = $newRow.Cells.Item(7)
$cell7.Range.Text = "Same as expected result."

# Cell 8: Pass/Fail (green highlight)
$cell8 = $newRow.Cells.Item(8)
$cell8.Range.Text = "Pass"
$cell8.Range.HighlightColorIndex = 4

# ---------------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the start of the new row.
# ---------------------------------------------------------------------------
$bmRange = $d.Range($cell1.Range.Start, $cell1.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
